$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 11666.333
$ws.Range("J6").Value = 11666.333
$ws.Range("L6").Value = 11666.333
$ws.Range("N6").Value = -12012.333

$ws.Range("H10").Value = 5950
$ws.Range("J10").Value = 5950
$ws.Range("L10").Value = 5950
$ws.Range("N10").Value = -6290

$ws.Range("H32").Value = 25006310
$ws.Range("I32").Value = 31252670
$ws.Range("J32").Value = 20874.625
$ws.Range("K32").Value = 31252670
$ws.Range("L32").Value = 20874.625
$ws.Range("M32").Value = -31252383
$ws.Range("N32").Value = -21448.625

$ws.Range("H102").Value = 3242.8572
$ws.Range("I102").Value = 2092.6667
$ws.Range("J102").Value = 4105.5
$ws.Range("K102").Value = 2092.6667
$ws.Range("L102").Value = 4105.5
$ws.Range("M102").Value = -470.6667000000002
$ws.Range("N102").Value = -7349.5

$ws.Range("H110").Value = 922.52
$ws.Range("I110").Value = 931.7917
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 931.7917
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 1113.2083
$ws.Range("N110").Value = -4790

$ws.Range("H132").Value = 2656.9678
$ws.Range("I132").Value = 2797.5
$ws.Range("J132").Value = 2507.0667
$ws.Range("K132").Value = 8392.5
$ws.Range("L132").Value = 7521.2001
$ws.Range("M132").Value = -5862.5
$ws.Range("N132").Value = -12581.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2859.16
$ws.Range("I105").Value = 1297.5
$ws.Range("J105").Value = 2994.9565
$ws.Range("K105").Value = 1297.5
$ws.Range("L105").Value = 2994.9565
$ws.Range("M105").Value = 449.5
$ws.Range("N105").Value = -6488.9565

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 83.29412000000001
$ws.Range("I7").Value = 64.72727
$ws.Range("J7").Value = 117.333336
$ws.Range("K7").Value = 64.72727
$ws.Range("L7").Value = 117.333336
$ws.Range("M7").Value = 48.27273
$ws.Range("N7").Value = -343.333336

$ws.Range("H12").Value = 10201
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 11586.857
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 11586.857
$ws.Range("M12").Value = -330
$ws.Range("N12").Value = -11926.857

$ws.Range("H42").Value = 8500
$ws.Range("I42").Value = 7000
$ws.Range("K42").Value = 7000
$ws.Range("M42").Value = -6407

$ws.Range("H105").Value = 716457.2
$ws.Range("I105").Value = 1112909.8
$ws.Range("J105").Value = 2842.6
$ws.Range("K105").Value = 1112909.8
$ws.Range("L105").Value = 2842.6
$ws.Range("M105").Value = -1111162.8
$ws.Range("N105").Value = -6336.6

$ws.Range("H107").Value = 1081.6364
$ws.Range("I107").Value = 995.25
$ws.Range("J107").Value = 1312
$ws.Range("K107").Value = 995.25
$ws.Range("L107").Value = 1312
$ws.Range("M107").Value = 924.75
$ws.Range("N107").Value = -5152

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1210.9814
$ws.Range("I68").Value = 1248
$ws.Range("J68").Value = 1190.8857
$ws.Range("K68").Value = 3744
$ws.Range("L68").Value = 3572.6571
$ws.Range("M68").Value = -2933
$ws.Range("N68").Value = -5194.6571

$ws.Range("H71").Value = 1210.9814
$ws.Range("I71").Value = 1248
$ws.Range("J71").Value = 1190.8857
$ws.Range("K71").Value = 11232
$ws.Range("L71").Value = 10717.9713
$ws.Range("M71").Value = -7176
$ws.Range("N71").Value = -18829.9713

$ws.Range("H122").Value = 13106.941
$ws.Range("I122").Value = 19721.727
$ws.Range("J122").Value = 979.8333
$ws.Range("K122").Value = 177495.543
$ws.Range("L122").Value = 8818.4997
$ws.Range("M122").Value = -175045.543
$ws.Range("N122").Value = -13718.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H45").Value = 30000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 30000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 30000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -31118

$ws.Range("H51").Value = 28666.334
$ws.Range("J51").Value = 28666.334
$ws.Range("L51").Value = 28666.334
$ws.Range("N51").Value = -29684.334

$ws.Range("H70").Value = 5004.2554
$ws.Range("I70").Value = 5147.93
$ws.Range("J70").Value = 3459.75
$ws.Range("K70").Value = 5147.93
$ws.Range("L70").Value = 3459.75
$ws.Range("M70").Value = -4877.93
$ws.Range("N70").Value = -3999.75

$ws.Range("H73").Value = 5004.2554
$ws.Range("I73").Value = 5147.93
$ws.Range("J73").Value = 3459.75
$ws.Range("K73").Value = 5147.93
$ws.Range("L73").Value = 3459.75
$ws.Range("M73").Value = -4211.93
$ws.Range("N73").Value = -5331.75

$ws.Range("H97").Value = 477311.94
$ws.Range("I97").Value = 556724.4399999999
$ws.Range("J97").Value = 837
$ws.Range("K97").Value = 556724.4399999999
$ws.Range("L97").Value = 837
$ws.Range("M97").Value = -556228.4399999999
$ws.Range("N97").Value = -1829

$ws.Range("H113").Value = 4726.1875
$ws.Range("I113").Value = 1400.9259
$ws.Range("J113").Value = 22682.6
$ws.Range("K113").Value = 1400.9259
$ws.Range("L113").Value = 22682.6
$ws.Range("M113").Value = 769.0741
$ws.Range("N113").Value = -27022.6

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 435100
$ws.Range("I33").Value = 650650
$ws.Range("K33").Value = 650650
$ws.Range("M33").Value = -650360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 414.2
$ws.Range("I107").Value = 424.46155
$ws.Range("J107").Value = 384.55554
$ws.Range("K107").Value = 1273.38465
$ws.Range("L107").Value = 1153.66662
$ws.Range("M107").Value = 646.61535
$ws.Range("N107").Value = -4993.66662

$ws.Range("H132").Value = 1564.7693
$ws.Range("I132").Value = 1127.7727
$ws.Range("J132").Value = 3968.25
$ws.Range("K132").Value = 3383.3181
$ws.Range("L132").Value = 11904.75
$ws.Range("M132").Value = -853.3181
$ws.Range("N132").Value = -16964.75
